$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 26, shifting existing rows 26-115 down to 27-116
$ws.Rows("26:26").Insert()

# Populate the new row 26 with values
$ws.Range("A26").Value = 9
$ws.Range("B26").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C26").Value = "Metropolitana"
$ws.Range("D26").Value = 44914
$ws.Range("E26").Value = 13
$ws.Range("F26").Value = 100114007
$ws.Range("G26").Value = "Jengibre"
$ws.Range("H26").Value = "Sin especificar"
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 520
$ws.Range("K26").Value = 13000
$ws.Range("L26").Value = 14000
$ws.Range("M26").Value = 13500
$ws.Range("N26").Value = "$/caja 13 kilos"
$ws.Range("O26").Value = "Perú"
$ws.Range("P26").Value = 1038
$ws.Range("Q26").Value = 13
$ws.Range("R26").Value = "Hortaliza"
